$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Record", "Balanço Geral", "Economia", "2025-04-08T13:12", "Positivo", "Vagas de emprego. Em Campos, são 391 vagas. Entre elas, para agente de saneamento, atendente, motorista de caminhão e auxiliar de cozinha. "),
    @("Record", "Balanço Geral", "Defesa Civil", "2025-04-08T13:19", "Neutro", "Tô na bronca. Moradores do bairro Ururaí reclamam de árvores em meio a fiação elétrica. *nota coberta*. Moradores disseram que já perderam eletrodomésticos. Concessionária disse que podas são de responsabilidade da prefeitura. Se estiver em contato com a rede elétrica, a empresa realiza o serviço. Uma equipe irá ao local. *com nota* da Defesa Civil, que vai enviar equipe para fazer vistoria das árvores. A responsabilidade pela poda é da Secretaria de serviços Públicos.    "),
    @("Record", "Balanço Geral", "Obras", "2025-04-08T13:21", "Negativo", "Tô na bronca. Moradores da Vila Romana reclamam de lama durante obras paralisadas. *nota coberta*. Situação teria se agravado porque calçamentos foram retirados. Em nota, prefeitura disse que as obras vêm sendo realizadas de forma simultânea. Prefeitura aguarda que concessionária conclua a parte dela. *com nota*"),
    @("Record", "Balanço Geral", "Infraestrutura", "2025-04-08T13:22", "Neutro", "BG mostrou, BG resolveu. Esgoto no Turfe é resolvido após ser exibido no Balanço Geral. Moradores reclamaram do esgoto da Rua Almeida Barbosa. problema durava há quase um mês. Exibido vídeo com a reclamação e com agradecimento da moradora dizendo que a Águas do Paraíba resolveu o problema. ")
)

$startRow = 81
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
}
